$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 6: the "1550免费分钟" business item description moves the
# "1550" number out of the short name and into the longer description.
$ws.Range("C6").Value = "免费分钟"
$ws.Range("D6").Value = "套餐不仅包含1550免费分钟，还有一条宽带可以使用"

# Move the active selection to H21, matching the saved cursor position.
$ws.Range("H21").Select()
